$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1
$ws.Range("A1").Value = 45436

# Update prices in the "COMUN" section
$ws.Range("D35").Value = 3666.444
$ws.Range("D36").Value = 4110
$ws.Range("D37").Value = 4110
$ws.Range("D38").Value = 4110
$ws.Range("D39").Value = 10440.086

# Update prices in the "CON TOPE" section
$ws.Range("D42").Value = 3953.09
$ws.Range("D43").Value = 4350
$ws.Range("D44").Value = 4350
$ws.Range("D45").Value = 4350
$ws.Range("D46").Value = 11676.446
